# Reproduce the meaningful, object-model-visible parts of the commit:
#  - "Run Data" (sheet1) data columns A:D narrowed slightly (13.5703125 -> 13.5 chars)
#  - "Run Info" (sheet2) column A narrowed the same way, column B widened/best-fit to
#    fit its longest label ("Run307_Sep21_15-03_Sector 4_Hole A.xlsx") at ~30.5 chars
#  - "Run Info" (sheet2) selection left on B16 (below the last used row) as the file
#    was saved
#
# (Cosmetic application/version stamps in the original diff -- fileVersion,
# mc:AlternateContent/x15 markup, theme font substitutions, calcPr concurrentCalc,
# x14ac:dyDescent rounding, indexed- vs auto- border color, and the "Run Data"
# sheetView topLeftCell scroll position -- are artifacts of the specific Excel
# build/OS that resaved the file and are not reachable/settable from the Excel
# COM object model exposed by this host, so they are intentionally left alone.)

$wb = $excel.ActiveWorkbook

# --- "Run Data" sheet -------------------------------------------------
$wsData = $wb.Worksheets.Item("Run Data")
$wsData.Columns("A:D").ColumnWidth = 12.7   # -> stored column width 13.5

# --- "Run Info" sheet --------------------------------------------------
$wsInfo = $wb.Worksheets.Item("Run Info")
$wsInfo.Columns("A:A").ColumnWidth = 12.7   # -> stored column width 13.5
$wsInfo.Columns("B:B").ColumnWidth = 29.6   # -> stored column width 30.5 (best-fit for Filename text)

# Leave the selection on "Run Info" (the active tab) at B16, matching the saved file.
$wsInfo.Activate()
$wsInfo.Range("B16").Select() | Out-Null
